{"js": "// Replace the \"Cyber Security\" themed document with the \"Chemistry\" themed\n// content, per the target revision: title, author, e-mail, body copy, and\n// the Summary section are all rewritten, and a trailing empty paragraph is\n// appended at the end of the document body.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Title\n// ---------------------------------------------------------------------\nconst titleHits = body.search(\"Cyber Security: The Digital Shield\", { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\ntitleHits.items[0].insertText(\n  \"Unveiling the Marvels of Chemistry: A Journey of Transformation and Discovery\",\n  Word.InsertLocation.replace\n);\n\n// ---------------------------------------------------------------------\n// Author name\n// ---------------------------------------------------------------------\nconst authorHits = body.search(\"Andrew Davenport\", { matchCase: true });\nauthorHits.load(\"items\");\nawait context.sync();\nauthorHits.items[0].insertText(\"Dr. Erika Martinez\", Word.InsertLocation.replace);\n\n// ---------------------------------------------------------------------\n// E-mail address\n// ---------------------------------------------------------------------\nconst emailHits = body.search(\"andrew.davenport@blackrock.net\", { matchCase: true });\nemailHits.load(\"items\");\nawait context.sync();\nemailHits.items[0].insertText(\"emartinez@highschooledu.org\", Word.InsertLocation.replace);\n\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Main body paragraph (the long essay with two manual line breaks)\n// ---------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst LF = \"\\v\"; // Word.js represents a <w:br/> with a vertical-tab char.\n\nconst essayParagraph = paragraphs.items[4];\nessayParagraph.load(\"text\");\nawait context.sync();\n\nconst newEssay =\n  \"In the vast tapestry of human knowledge, chemistry stands as a vibrant and transformative thread, weaving together the elements of the universe to reveal the wonders of matter and its intricate interactions.\" +\n  \" It is the language of nature's symphony, guiding us through the profound mysteries of chemical reactions, the compositions of substances, and the fundamental principles that govern the behaviour of matter.\" +\n  \" As we embark on this enthralling adventure into the world of chemistry, let us unravel the enigmatic secrets of this fascinating science, unveiling its profound impact on our lives.\" +\n  LF + LF +\n  \"In the realm of chemistry, we witness the ceaseless dance of atoms and molecules, orchestrated by the enigmatic forces of attraction and repulsion.\" +\n  \" These fundamental particles engage in a ceaseless choreography, forming compounds and rearranging their structures, giving rise to the mesmerizing diversity of substances that make up our world.\" +\n  \" From the air we breathe and the water we drink, to the food we consume and the medicines that heal us, chemistry lies at the heart of all life.\" +\n  \" It is the driving force behind the symphony of life, governing the intricate interactions between organisms in the intricate web of ecosystems.\" +\n  LF + LF +\n  \"Through the lens of chemistry, we can unravel the complexities of chemical reactions, understanding how substances transform from one state to another, releasing energy or undergoing profound changes in their properties.\" +\n  \" We uncover the secrets of catalysts, molecules that accelerate these transformations, allowing us to harness nature's power to create new substances and materials.\" +\n  \" Chemistry empowers us with the ability to synthesize drugs that combat diseases, develop materials with extraordinary properties, and create sustainable energy sources, all of which have the potential to shape a better future for humankind.\";\n\nessayParagraph.insertText(newEssay, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// \"Summary\" heading stays the same text; only the body below changes.\n// ---------------------------------------------------------------------\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst summaryParagraph = paragraphs.items[6];\nsummaryParagraph.load(\"text\");\nawait context.sync();\n\nconst newSummary =\n  \"Chemistry is the vibrant language of nature, revealing the secrets of matter and its intricate interactions.\" +\n  \" From the symphony of atomic dances to the transformative power of chemical reactions, chemistry plays a pivotal role in shaping our world.\" +\n  \" It holds the key to understanding the composition and behaviour of substances, leading to advancements in medicine, technology, and sustainable energy.\" +\n  \" Chemistry empowers us to comprehend and manipulate the world around us, fostering progress and innovation while deepening our appreciation for the marvels of nature.\";\n\nsummaryParagraph.insertText(newSummary, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Append a new, empty trailing paragraph at the end of the document body.\n// ---------------------------------------------------------------------\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Replace the \"Cyber Security\" themed document with the \"Chemistry\" themed\n# content, per the target revision: title, author, e-mail, body copy, and\n# the Summary section are all rewritten, and a trailing empty paragraph is\n# appended at the end of the document body.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Title\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Cyber Security: The Digital Shield\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Unveiling the Marvels of Chemistry: A Journey of Transformation and Discovery\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# ---------------------------------------------------------------------\n# Author name\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Andrew Davenport\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Dr. Erika Martinez\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# ---------------------------------------------------------------------\n# E-mail address\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"andrew.davenport@blackrock.net\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"emartinez@highschooledu.org\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# ---------------------------------------------------------------------\n# Main body paragraph (the long essay with two manual line breaks)\n# ---------------------------------------------------------------------\n$LF = [char]11   # Word.js / COM manual line break <w:br/>\n\n$essayParagraph = $d.Paragraphs(5)\n$essayRange = $essayParagraph.Range\n$essayRange.MoveEnd(1, -1)   # keep the trailing paragraph mark untouched\n\n$newEssay = (\n  \"In the vast tapestry of human knowledge, chemistry stands as a vibrant and transformative thread, weaving together the elements of the universe to reveal the wonders of matter and its intricate interactions.\" +\n  \" It is the language of nature's symphony, guiding us through the profound mysteries of chemical reactions, the compositions of substances, and the fundamental principles that govern the behaviour of matter.\" +\n  \" As we embark on this enthralling adventure into the world of chemistry, let us unravel the enigmatic secrets of this fascinating science, unveiling its profound impact on our lives.\" +\n  $LF + $LF +\n  \"In the realm of chemistry, we witness the ceaseless dance of atoms and molecules, orchestrated by the enigmatic forces of attraction and repulsion.\" +\n  \" These fundamental particles engage in a ceaseless choreography, forming compounds and rearranging their structures, giving rise to the mesmerizing diversity of substances that make up our world.\" +\n  \" From the air we breathe and the water we drink, to the food we consume and the medicines that heal us, chemistry lies at the heart of all life.\" +\n  \" It is the driving force behind the symphony of life, governing the intricate interactions between organisms in the intricate web of ecosystems.\" +\n  $LF + $LF +\n  \"Through the lens of chemistry, we can unravel the complexities of chemical reactions, understanding how substances transform from one state to another, releasing energy or undergoing profound changes in their properties.\" +\n  \" We uncover the secrets of catalysts, molecules that accelerate these transformations, allowing us to harness nature's power to create new substances and materials.\" +\n  \" Chemistry empowers us with the ability to synthesize drugs that combat diseases, develop materials with extraordinary properties, and create sustainable energy sources, all of which have the potential to shape a better future for humankind.\"\n)\n\n$essayRange.Text = $newEssay\n\n# ---------------------------------------------------------------------\n# \"Summary\" heading stays the same text; only the body below changes.\n# ---------------------------------------------------------------------\n$summaryParagraph = $d.Paragraphs(7)\n$summaryRange = $summaryParagraph.Range\n$summaryRange.MoveEnd(1, -1)   # keep the trailing paragraph mark untouched\n\n$newSummary = (\n  \"Chemistry is the vibrant language of nature, revealing the secrets of matter and its intricate interactions.\" +\n  \" From the symphony of atomic dances to the transformative power of chemical reactions, chemistry plays a pivotal role in shaping our world.\" +\n  \" It holds the key to understanding the composition and behaviour of substances, leading to advancements in medicine, technology, and sustainable energy.\" +\n  \" Chemistry empowers us to comprehend and manipulate the world around us, fostering progress and innovation while deepening our appreciation for the marvels of nature.\"\n)\n\n$summaryRange.Text = $newSummary\n\n# ---------------------------------------------------------------------\n# Append a new, empty trailing paragraph at the end of the document body.\n# ---------------------------------------------------------------------\n$endRange = $d.Content\n$endRange.Collapse(0)   # wdCollapseEnd\n$endRange.InsertParagraphAfter()\n"}
